# Fruta / hortaliza, semanal
#
# A new weekly price record for "Arveja Verde" (Macroferia Regional de Talca)
# is inserted as row 50, pushing the existing rows 50-68 down to 51-69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 50 (shifts rows 50:68 -> 51:69).
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new weekly record.
$ws.Cells.Item(50, 1).Value = 5
$ws.Cells.Item(50, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(50, 3).Value = "Maule"
$ws.Cells.Item(50, 4).Value = 44524
$ws.Cells.Item(50, 5).Value = 7
$ws.Cells.Item(50, 6).Value = 100112022
$ws.Cells.Item(50, 7).Value = "Arveja Verde"
$ws.Cells.Item(50, 8).Value = "Sin especificar"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 500
$ws.Cells.Item(50, 11).Value = 14000
$ws.Cells.Item(50, 12).Value = 14000
$ws.Cells.Item(50, 13).Value = 14000
$ws.Cells.Item(50, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(50, 15).Value = "Región del Maule"
$ws.Cells.Item(50, 16).Value = 560
$ws.Cells.Item(50, 17).Value = 25
$ws.Cells.Item(50, 18).Value = "Hortaliza"
